# Fruta / hortaliza, semanal
# The weekly refresh re-shuffles which data row (2-12) carries which
# Fecha/Volumen/Precio set. Columns D (Fecha), L (Calidad), M (Volumen),
# N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado) and
# S (Precio $/Kg) are re-mapped across rows 2..12 according to the table
# below; every other column stays untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together as a "record" when rows are re-ordered.
$cols = @("D", "L", "M", "N", "O", "P", "S")

$firstRow = 2
$lastRow = 12

# Snapshot the current (before) values for the columns that change, since
# the remapping below is a permutation with cycles and must not read back
# values that were already overwritten earlier in the loop.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Maps each destination row to the source row whose record it now holds.
$rowMap = @{
    2  = 5
    3  = 4
    4  = 10
    5  = 7
    6  = 8
    7  = 2
    8  = 9
    9  = 11
    10 = 12
    11 = 3
    12 = 6
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $srcVals[$c]
    }
}
